$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.985.74"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.345.22"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.13"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.31"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").Value = "3.336.03"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.626"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.55"
$ws.Range("E12").Value = "  -4.82%  "
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.08"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "3.882.40"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.03"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "3.354.84"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "64.988.41"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.65"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.984"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "474.80"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.89"
$ws.Range("E23").Value = "  -5.65%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.04"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.00"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.45"
$ws.Range("E28").Value = "  -4.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.57"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.87"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.44"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.29"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "61.69"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "571.00"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.13"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.368"
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").Value = "0.0₃0727"
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").Value = "3.074.74"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0410"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.17"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.11"
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -0.81%  "
